$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "20.510.71"
Set-TextValue "E2" "  +2.86%  "
Set-TextValue "D3" "1.469.78"
Set-TextValue "E3" "  +3.51%  "
Set-TextValue "E4" "  +0.33%  "
Set-TextValue "D5" "0.9563"
Set-TextValue "E5" "  -4.67%  "
Set-TextValue "D6" "281.83"
Set-TextValue "E6" "  +3.20%  "
Set-TextValue "D7" "0.3697"
Set-TextValue "E7" "  -1.30%  "
Set-TextValue "D8" "0.3180"
Set-TextValue "E8" "  +2.84%  "
Set-TextValue "D9" "41.82"
Set-TextValue "E9" "  +4.65%  "
Set-TextValue "D10" "1.055"
Set-TextValue "E10" "  +4.35%  "
Set-TextValue "E11" "  +1.34%  "
Set-TextValue "E12" "  -0.10%  "
Set-TextValue "D13" "5.610"
Set-TextValue "E13" "  +4.20%  "
Set-TextValue "D14" "18.21"
Set-TextValue "E14" "  +6.74%  "
Set-TextValue "B15" "Chainlink"
Set-TextValue "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "6.261"
Set-TextValue "E15" "  +1.65%  "
Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "1.473.91"
Set-TextValue "E16" "  +3.31%  "
Set-TextValue "E17" "  +2.89%  "
Set-TextValue "D18" "0.05693"
Set-TextValue "E18" "  -2.46%  "
Set-TextValue "B19" "Dai"
Set-TextValue "C19" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D19" "0.9565"
Set-TextValue "E19" "  -4.65%  "
Set-TextValue "B20" "Litecoin"
Set-TextValue "C20" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D20" "72.20"
Set-TextValue "E20" "  -3.82%  "
Set-TextValue "D21" "5.681"
Set-TextValue "E21" "  +0.40%  "
Set-TextValue "D22" "14.72"
Set-TextValue "E22" "  +1.96%  "
Set-TextValue "D23" "11.21"
Set-TextValue "E23" "  +1.99%  "
Set-TextValue "D24" "2.268"
Set-TextValue "E24" "  -3.06%  "
Set-TextValue "D25" "20.676.68"
Set-TextValue "E25" "  +3.54%  "
Set-TextValue "D26" "2.294"
Set-TextValue "E26" "  +0.54%  "
Set-TextValue "D27" "137.82"
Set-TextValue "E27" "  -0.92%  "
Set-TextValue "D28" "17.55"
Set-TextValue "E28" "  +4.11%  "
Set-TextValue "D29" "1.637.92"
Set-TextValue "E29" "  +3.18%  "
Set-TextValue "D30" "113.83"
Set-TextValue "E30" "  +4.17%  "
Set-TextValue "D31" "3.952"
Set-TextValue "E31" "  +2.14%  "
Set-TextValue "D32" "5.320"
Set-TextValue "E32" "  -1.67%  "
Set-TextValue "D33" "0.8337"
Set-TextValue "E33" "  -6.37%  "
Set-TextValue "D34" "1.623"
Set-TextValue "E34" "  +27.43%  "
Set-TextValue "D35" "0.07840"
Set-TextValue "E35" "  +0.76%  "
Set-TextValue "D36" "0.06045"
Set-TextValue "E36" "  +6.37%  "
Set-TextValue "D37" "4.913"
Set-TextValue "E37" "  +3.49%  "
Set-TextValue "E38" "  +3.17%  "
Set-TextValue "D39" "10.60"
Set-TextValue "E39" "  -5.65%  "
Set-TextValue "B40" "Frax"
Set-TextValue "C40" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D40" "0.9696"
Set-TextValue "E40" "  -3.30%  "
Set-TextValue "B41" "TrustWalletToken"
Set-TextValue "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D41" "1.122"
Set-TextValue "E41" "  +2.37%  "
Set-TextValue "E42" "  -1.15%  "
Set-TextValue "D43" "7.352"
Set-TextValue "E43" "  -12.16%  "
Set-TextValue "D44" "0.5419"
Set-TextValue "E44" "  +1.92%  "
Set-TextValue "D45" "12.50"
Set-TextValue "E45" "  +1.74%  "
Set-TextValue "D46" "3.593"
Set-TextValue "E46" "  +1.69%  "
Set-TextValue "D47" "122.16"
Set-TextValue "E47" "  +11.38%  "
Set-TextValue "D48" "0.5331"
Set-TextValue "E48" "  +4.14%  "
Set-TextValue "D49" "1.831"
Set-TextValue "E49" "  +2.45%  "
Set-TextValue "D50" "0.06436"
Set-TextValue "E50" "  +4.29%  "
Set-TextValue "D51" "1.051"
Set-TextValue "E51" "  +0.09%  "
